$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new header columns (D, E, F) on row 1: ORG_TOW_IDENOLD,
# ORG_TOW_IDENNEW, ORG_TOW_STATUS, left-aligned like the existing A1:C1
# headers (style index 1 in the original file == horizontal="left").
$ws.Range("D1:F1").HorizontalAlignment = -4131

$ws.Range("D1").Value = "ORG_TOW_IDENOLD"
$ws.Range("E1").Value = "ORG_TOW_IDENNEW"
$ws.Range("F1").Value = "ORG_TOW_STATUS"

# The new columns' data rows (2:19) are left empty, matching the source.

# Restore the active-cell selection recorded in the saved file's view state.
$ws.Range("F4").Select()
